# Append 8 new species-observation records (rows 9-16) to the "Artfynd" sheet,
# extending the used range from A1:AY8 to A1:AY16.
#
# Values that Excel would otherwise auto-convert (pure digit strings in the
# "Antal" column, ISO dates/times, and explicit empty strings) are written
# with a leading single-quote so they stay text cells, matching the source
# data's inline-string typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 112154275
$ws.Range("B9").Value = 44322
$ws.Range("C9").Value = 'Ovaliderad'
$ws.Range("D9").Value = 'NT'
$ws.Range("E9").Value = 102366
$ws.Range("F9").Value = 'Ängsmetallvinge'
$ws.Range("G9").Value = 'Adscita statices'
$ws.Range("H9").Value = '(Linnaeus, 1758)'
$ws.Range("I9").Value = "'5"
$ws.Range("J9").Value = 'ex.'
$ws.Range("M9").Value = 'vilande'
$ws.Range("P9").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q9").Value = 442616.138687243
$ws.Range("R9").Value = 6204441.08982533
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = 'Skåne'
$ws.Range("U9").Value = 'Kristianstad'
$ws.Range("V9").Value = 'Skåne'
$ws.Range("W9").Value = 'Vä'
$ws.Range("Y9").Value = "'2013-07-09"
$ws.Range("Z9").Value = "'00:00"
$ws.Range("AA9").Value = "'2013-07-09"
$ws.Range("AB9").Value = "'00:00"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AI9").Value = 'på igenväxande sandhed'
$ws.Range("AO9").Value = 'på blmr av åkervädd m fl'
$ws.Range("AT9").Value = "'"
$ws.Range("AW9").Value = 'Nils Otto Nilsson'
$ws.Range("AX9").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY9").Value = 'Krst NV-program 2013'

$ws.Range("A10").Value = 112154273
$ws.Range("B10").Value = 44322
$ws.Range("C10").Value = 'Ovaliderad'
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 102366
$ws.Range("F10").Value = 'Ängsmetallvinge'
$ws.Range("G10").Value = 'Adscita statices'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("I10").Value = "'4"
$ws.Range("J10").Value = 'ex.'
$ws.Range("M10").Value = 'vilande'
$ws.Range("P10").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q10").Value = 442541.7120545401
$ws.Range("R10").Value = 6204451.031370129
$ws.Range("S10").Value = 10
$ws.Range("T10").Value = 'Skåne'
$ws.Range("U10").Value = 'Kristianstad'
$ws.Range("V10").Value = 'Skåne'
$ws.Range("W10").Value = 'Vä'
$ws.Range("Y10").Value = "'2013-07-09"
$ws.Range("Z10").Value = "'00:00"
$ws.Range("AA10").Value = "'2013-07-09"
$ws.Range("AB10").Value = "'00:00"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AG10").Value = $false
$ws.Range("AI10").Value = 'på igenväxande sandhed'
$ws.Range("AO10").Value = 'på blmr av åkervädd m fl'
$ws.Range("AT10").Value = "'"
$ws.Range("AW10").Value = 'Nils Otto Nilsson'
$ws.Range("AX10").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY10").Value = 'Krst NV-program 2013'

$ws.Range("A11").Value = 112154281
$ws.Range("B11").Value = 42578
$ws.Range("C11").Value = 'Ovaliderad'
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 100770
$ws.Range("F11").Value = 'Mindre blåvinge'
$ws.Range("G11").Value = 'Cupido minimus'
$ws.Range("H11").Value = '(Fuessly, 1775)'
$ws.Range("I11").Value = "'1"
$ws.Range("J11").Value = 'ex.'
$ws.Range("M11").Value = 'friflygande'
$ws.Range("P11").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q11").Value = 442664.1890363992
$ws.Range("R11").Value = 6204260.315617888
$ws.Range("S11").Value = 10
$ws.Range("T11").Value = 'Skåne'
$ws.Range("U11").Value = 'Kristianstad'
$ws.Range("V11").Value = 'Skåne'
$ws.Range("W11").Value = 'Vä'
$ws.Range("Y11").Value = "'2013-07-09"
$ws.Range("Z11").Value = "'00:00"
$ws.Range("AA11").Value = "'2013-07-09"
$ws.Range("AB11").Value = "'00:00"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AI11").Value = 'på igenväxande sandhed'
$ws.Range("AT11").Value = "'"
$ws.Range("AW11").Value = 'Nils Otto Nilsson'
$ws.Range("AX11").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY11").Value = 'Krst NV-program 2013'

$ws.Range("A12").Value = 112154283
$ws.Range("B12").Value = 44331
$ws.Range("C12").Value = 'Ovaliderad'
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 201164
$ws.Range("F12").Value = 'Sexfläckig bastardsvärmare'
$ws.Range("G12").Value = 'Zygaena filipendulae'
$ws.Range("H12").Value = '(Linnaeus, 1758)'
$ws.Range("I12").Value = "'1"
$ws.Range("J12").Value = 'ex.'
$ws.Range("M12").Value = 'födosökande'
$ws.Range("P12").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q12").Value = 442664.1890363992
$ws.Range("R12").Value = 6204260.315617888
$ws.Range("S12").Value = 10
$ws.Range("T12").Value = 'Skåne'
$ws.Range("U12").Value = 'Kristianstad'
$ws.Range("V12").Value = 'Skåne'
$ws.Range("W12").Value = 'Vä'
$ws.Range("Y12").Value = "'2013-07-09"
$ws.Range("Z12").Value = "'00:00"
$ws.Range("AA12").Value = "'2013-07-09"
$ws.Range("AB12").Value = "'00:00"
$ws.Range("AD12").Value = $false
$ws.Range("AE12").Value = $false
$ws.Range("AG12").Value = $false
$ws.Range("AI12").Value = 'på igenväxande sandhed'
$ws.Range("AO12").Value = 'på blmr av åkervädd'
$ws.Range("AT12").Value = "'"
$ws.Range("AW12").Value = 'Nils Otto Nilsson'
$ws.Range("AX12").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY12").Value = 'Krst NV-program 2013'

$ws.Range("A13").Value = 112154272
$ws.Range("B13").Value = 44322
$ws.Range("C13").Value = 'Ovaliderad'
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 102366
$ws.Range("F13").Value = 'Ängsmetallvinge'
$ws.Range("G13").Value = 'Adscita statices'
$ws.Range("H13").Value = '(Linnaeus, 1758)'
$ws.Range("I13").Value = "'3"
$ws.Range("J13").Value = 'ex.'
$ws.Range("M13").Value = 'vilande'
$ws.Range("P13").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q13").Value = 442480.1045430943
$ws.Range("R13").Value = 6204371.296283903
$ws.Range("S13").Value = 10
$ws.Range("T13").Value = 'Skåne'
$ws.Range("U13").Value = 'Kristianstad'
$ws.Range("V13").Value = 'Skåne'
$ws.Range("W13").Value = 'Vä'
$ws.Range("Y13").Value = "'2013-07-09"
$ws.Range("Z13").Value = "'00:00"
$ws.Range("AA13").Value = "'2013-07-09"
$ws.Range("AB13").Value = "'00:00"
$ws.Range("AD13").Value = $false
$ws.Range("AE13").Value = $false
$ws.Range("AG13").Value = $false
$ws.Range("AI13").Value = 'på igenväxande sandhed'
$ws.Range("AO13").Value = 'på blmr av åkervädd m fl'
$ws.Range("AT13").Value = "'"
$ws.Range("AW13").Value = 'Nils Otto Nilsson'
$ws.Range("AX13").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY13").Value = 'Krst NV-program 2013'

$ws.Range("A14").Value = 112154276
$ws.Range("B14").Value = 39449
$ws.Range("C14").Value = 'Ovaliderad'
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 102471
$ws.Range("F14").Value = 'Åkerväddsantennmal'
$ws.Range("G14").Value = 'Nemophora metallica'
$ws.Range("H14").Value = '(Poda, 1761)'
$ws.Range("I14").Value = "'1"
$ws.Range("J14").Value = 'ex.'
$ws.Range("M14").Value = 'vilande'
$ws.Range("P14").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q14").Value = 442616.138687243
$ws.Range("R14").Value = 6204441.08982533
$ws.Range("S14").Value = 10
$ws.Range("T14").Value = 'Skåne'
$ws.Range("U14").Value = 'Kristianstad'
$ws.Range("V14").Value = 'Skåne'
$ws.Range("W14").Value = 'Vä'
$ws.Range("Y14").Value = "'2013-07-09"
$ws.Range("Z14").Value = "'00:00"
$ws.Range("AA14").Value = "'2013-07-09"
$ws.Range("AB14").Value = "'00:00"
$ws.Range("AD14").Value = $false
$ws.Range("AE14").Value = $false
$ws.Range("AG14").Value = $false
$ws.Range("AI14").Value = 'på igenväxande sandhed'
$ws.Range("AO14").Value = 'på blmr av åkervädd'
$ws.Range("AT14").Value = "'"
$ws.Range("AW14").Value = 'Nils Otto Nilsson'
$ws.Range("AX14").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY14").Value = 'Krst NV-program 2013'

$ws.Range("A15").Value = 112154282
$ws.Range("B15").Value = 44322
$ws.Range("C15").Value = 'Ovaliderad'
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 102366
$ws.Range("F15").Value = 'Ängsmetallvinge'
$ws.Range("G15").Value = 'Adscita statices'
$ws.Range("H15").Value = '(Linnaeus, 1758)'
$ws.Range("I15").Value = "'4"
$ws.Range("J15").Value = 'ex.'
$ws.Range("M15").Value = 'vilande'
$ws.Range("P15").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q15").Value = 442664.1890363992
$ws.Range("R15").Value = 6204260.315617888
$ws.Range("S15").Value = 10
$ws.Range("T15").Value = 'Skåne'
$ws.Range("U15").Value = 'Kristianstad'
$ws.Range("V15").Value = 'Skåne'
$ws.Range("W15").Value = 'Vä'
$ws.Range("Y15").Value = "'2013-07-09"
$ws.Range("Z15").Value = "'00:00"
$ws.Range("AA15").Value = "'2013-07-09"
$ws.Range("AB15").Value = "'00:00"
$ws.Range("AD15").Value = $false
$ws.Range("AE15").Value = $false
$ws.Range("AG15").Value = $false
$ws.Range("AI15").Value = 'på igenväxande sandhed'
$ws.Range("AO15").Value = 'på blmr av åkervädd'
$ws.Range("AT15").Value = "'"
$ws.Range("AW15").Value = 'Nils Otto Nilsson'
$ws.Range("AX15").Value = 'Nils Otto Nilsson, Mats Karlsson'
$ws.Range("AY15").Value = 'Krst NV-program 2013'

$ws.Range("A16").Value = 112145588
$ws.Range("B16").Value = 42578
$ws.Range("C16").Value = 'Ovaliderad'
$ws.Range("D16").Value = 'NT'
$ws.Range("E16").Value = 100770
$ws.Range("F16").Value = 'Mindre blåvinge'
$ws.Range("G16").Value = 'Cupido minimus'
$ws.Range("H16").Value = '(Fuessly, 1775)'
$ws.Range("I16").Value = "'1"
$ws.Range("J16").Value = 'ex.'
$ws.Range("K16").Value = 'imago/adult'
$ws.Range("L16").Value = 'hane'
$ws.Range("M16").Value = 'födosökande'
$ws.Range("P16").Value = 'Vä, delomr 22, 700 m NNO Sigridslund, Sk'
$ws.Range("Q16").Value = 442602.7188004656
$ws.Range("R16").Value = 6204401.54979612
$ws.Range("S16").Value = 10
$ws.Range("T16").Value = 'Skåne'
$ws.Range("U16").Value = 'Kristianstad'
$ws.Range("V16").Value = 'Skåne'
$ws.Range("W16").Value = 'Vä'
$ws.Range("Y16").Value = "'2013-06-12"
$ws.Range("Z16").Value = "'00:00"
$ws.Range("AA16").Value = "'2013-06-12"
$ws.Range("AB16").Value = "'00:00"
$ws.Range("AD16").Value = $false
$ws.Range("AE16").Value = $false
$ws.Range("AG16").Value = $false
$ws.Range("AI16").Value = 'i igenväxande hed'
$ws.Range("AO16").Value = 'på praktveronika'
$ws.Range("AT16").Value = "'"
$ws.Range("AW16").Value = 'Nils Otto Nilsson'
$ws.Range("AX16").Value = 'Nils Otto Nilsson'
$ws.Range("AY16").Value = 'Krst NV-program 2013'
